# edit clickEvent(), edit save() +add saving tree
$wb = $excel.ActiveWorkbook

# --- "info" sheet: fill in A1:C1 (D1 already holds "3") ---
$infoSheet = $wb.Worksheets.Item("info")
$infoSheet.Range("A1").Value = "we"
$infoSheet.Range("B1").Value = "rwe"
$infoSheet.Range("C1").Value = "wer"

# --- "items" sheet: add the saving tree rows ---
$itemsSheet = $wb.Worksheets.Item("items")
$itemsSheet.Range("A1").Value = "접시140"
$itemsSheet.Range("B1").Value = "개"
$itemsSheet.Range("C1").Value = 1000
$itemsSheet.Range("D1").Value = 1
$itemsSheet.Range("E1").Value = 1000

$itemsSheet.Range("A2").Value = "부의금가방"
$itemsSheet.Range("B2").Value = "줄"
$itemsSheet.Range("C2").Value = 5000
$itemsSheet.Range("D2").Value = 1
$itemsSheet.Range("E2").Value = 5000
